$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (GitHub Actions scheduled update).
# Values are set with a leading text-qualifier apostrophe where the new
# value looks numeric, so Excel stores it as text (matching the column's
# existing inline-string / text formatting) instead of silently coercing
# it to a number and dropping meaningful trailing zeros (e.g. "0.1340").

$ws.Range('D2').Value = '29.361.05'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.839.84'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('D4').Value = '''0.9998'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''238.96'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').Value = '''0.6259'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').Value = '''0.07373'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').Value = '''0.2887'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('D10').Value = '''24.77'
$ws.Range('E10').Value = '  +1.56%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '1.831.51'
$ws.Range('E12').Value = '  -0.57%  '
$ws.Range('D13').Value = '''4.958'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').Value = '''0.6719'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = '''0.00001019'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').Value = '''81.71'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '''6.268'
$ws.Range('E17').Value = '  +1.47%  '
$ws.Range('D18').Value = '29.333.69'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Value = '''233.67'
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('D20').Value = '''12.31'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').Value = '''7.281'
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '''157.43'
$ws.Range('E24').Value = '  -0.67%  '
$ws.Range('D25').Value = '''8.479'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').Value = '''0.1340'
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('D27').Value = '''17.29'
$ws.Range('E27').Value = '  -1.10%  '
$ws.Range('D28').Value = '''0.07215'
$ws.Range('E28').Value = '  +12.53%  '
$ws.Range('D29').Value = '''1.489'
$ws.Range('E29').Value = '  +5.08%  '
$ws.Range('D30').Value = '''1.473'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').Value = '''4.032'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('D32').Value = '''4.020'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').Value = '''1.814'
$ws.Range('D34').Value = '''1.147'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('D35').Value = '''0.6976'
$ws.Range('E35').Value = '  +0.64%  '
$ws.Range('D36').Value = '''2.580'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '''0.01827'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('D38').Value = '''2.804'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').Value = '1.230.66'
$ws.Range('E39').Value = '  -2.23%  '
$ws.Range('D40').Value = '''6.734'
$ws.Range('E40').Value = '  -0.51%  '
$ws.Range('D41').Value = '''0.9429'
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').Value = '1.990.84'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('D44').Value = '''100.97'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').Value = '''65.25'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('E46').Value = '  +5.02%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value = '''6.943'
$ws.Range('E47').Value = '  -1.62%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.690'
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('D49').Value = '''8.887'
$ws.Range('E49').Value = '  -1.20%  '
$ws.Range('D50').Value = '''0.3876'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').Value = '''0.1124'
$ws.Range('E51').Value = '  -2.67%  '
